# Re-load the "quadratic-svm-score" sheet with an updated copy of the
# upstream ful-path.csv scores. The refreshed CSV import re-applies the
# text number format to the label column / header row (hence the style
# churn you see on A1:C1 and A2:A8) and brings in new decision-function
# scores for the "1-o__Chitinivibrionales" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns/headers are re-imported as text, same as before, just re-applied.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A8").NumberFormat = "@"

# Updated quadratic-svm-score values from the refreshed ful-path.csv.
$ws.Range("B2").Value = -2204157.7456578161
$ws.Range("B3").Value = -1609728.014645484
$ws.Range("B4").Value = -731123.63370493101
$ws.Range("B5").Value = -711149.02978820261
$ws.Range("B6").Value = 653556.79410539288
